$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211277842521667
$ws.Range("B1").Value = 2.613821029663086
$ws.Range("D1").Value = 2.172142505645752
$ws.Range("E1").Value = 1.16086757183075
